$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update "Förändrad" column C for rows 2-14: 46078 -> 46079
for ($r = 2; $r -le 14; $r++) {
    $ws.Cells.Item($r, 3).Value = 46079
}

# 2) Reorder rows 7-11 data (columns A, B, G) - cyclic shift
#    before: row7=62884-2021, row8=25634-2025, row9=28266-2025, row10=25015-2023, row11=19922-2025
#    after:  row7=19922-2025, row8=25015-2023, row9=62884-2021, row10=25634-2025, row11=28266-2025
$rows_7_11 = @(
    @{ A = "A 62884-2021"; B = 44504;              G = 0.8 },
    @{ A = "A 25634-2025"; B = 45803.59570601852;  G = 6 },
    @{ A = "A 28266-2025"; B = 45818.56381944445;  G = 1.9 },
    @{ A = "A 25015-2023"; B = 45085.6989699074;   G = 1.8 },
    @{ A = "A 19922-2025"; B = 45771.63034722222;  G = 10.1 }
)

# New order of source rows (0-indexed into $rows_7_11) for target rows 7..11
$order_7_11 = @(4, 3, 0, 1, 2)

for ($i = 0; $i -lt 5; $i++) {
    $targetRow = 7 + $i
    $src = $rows_7_11[$order_7_11[$i]]
    $ws.Cells.Item($targetRow, 1).Value = $src.A
    $ws.Cells.Item($targetRow, 2).Value = $src.B
    $ws.Cells.Item($targetRow, 7).Value = $src.G
}

# 3) Reorder rows 12-14 data (columns A, B, G) - cyclic shift
#    before: row12=60024-2025, row13=3402-2026, row14=14271-2021
#    after:  row12=14271-2021, row13=60024-2025, row14=3402-2026
$rows_12_14 = @(
    @{ A = "A 60024-2025"; B = 45992;              G = 1.1 },
    @{ A = "A 3402-2026";  B = 46042.39047453704;  G = 5.5 },
    @{ A = "A 14271-2021"; B = 44278;              G = 6.7 }
)

$order_12_14 = @(2, 0, 1)

for ($i = 0; $i -lt 3; $i++) {
    $targetRow = 12 + $i
    $src = $rows_12_14[$order_12_14[$i]]
    $ws.Cells.Item($targetRow, 1).Value = $src.A
    $ws.Cells.Item($targetRow, 2).Value = $src.B
    $ws.Cells.Item($targetRow, 7).Value = $src.G
}
